$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new data row (row 2) -------------------------------------------------
$ws.Range("A2").Value = "MCH110"
$ws.Range("C2").Value = "1. NGOMNGA, 5. STEVE BIKO ARTICLES WETHU U STEVE BANTU BIKO SISAMKHUMBULA NAMANJE, 2. NEWSPAPER CLIPPINGS STEVE BIKO, 3. STEVE BIKO DIES IN DETENTION, 4. TWO DIE AS YOUTH GOES ON RAMPAGE"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

# D2 and H2 stay blank but still pick up the row's formatting (see loop below)

# --- Apply the row's font (Calibri 10, automatic/theme text color) ----------------
foreach ($addr in @("A2", "C2", "D2", "E2", "F2", "G2", "H2")) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.Font.ThemeColor = 1
}

# --- Re-establish the frozen header row / selection --------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A5").Select()
